$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(3).Insert()
$ws.Range("C1").Value = "Chord"
for ($r = 2; $r -le 168; $r++) {
    $b = $ws.Cells.Item($r, 2).Text
    if ($b -eq "-") {
        $ws.Cells.Item($r, 3).Value = "-"
    } else {
        $ws.Cells.Item($r, 3).Value = $false
    }
}
